$wb = $excel.ActiveWorkbook

# Move the current selection on the 3rd tab and make that tab the active one
# (becomes bookViews/workbookView activeTab="2").
[void]$wb.Worksheets.Item("Project3").Range("B30").Select()

# Shorten the project sheet names; "Resources" is left unchanged.
$wb.Worksheets.Item("Project1").Name = "P1"
$wb.Worksheets.Item("Project2").Name = "P2"
$wb.Worksheets.Item("Project3").Name = "P3"
